$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hobbies answer (C3) - new delimiter style and values
$ws.Range("C3").Value = "Football||Coding||Cricket||Video Games||TV Series||Travelling"

# Update the "Tell us something about yourself" answer (C4) with new lorem text
$ws.Range("C4").Value = "Lorem Ipsum is simply dummy text of the printing and typesetting industry. Lorem Ipsum has been the industry's standard dummy text ever since the 1500s, when an unknown printer took a galley of type and scrambled it to make a type specimen book. It has survived not only five centuries, but also the leap into electronic typesetting, remaining essentially unchanged. It was popularised in the 1960s with the release of Letraset sheets containing Lorem Ipsum passages, and more recently with desktop publishing software like Aldus PageMaker including versions of Lorem Ipsum."

# Update the rating value (C5) - keep it stored as text (shared string), not a number
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "3.5"
$ws.Range("C5").ClearFormats()

# Remove row 6 entirely (the "Describe the above Image" question/answer row)
$ws.Range("A6:C6").EntireRow.Delete()

# Adjust column B width to match new content (best-fit shrink after row 6 removal)
$ws.Columns.Item(2).ColumnWidth = 31.5
